$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 212 and 213 (existing rows 212..264 shift down to 214..266)
$ws.Range("A212:A213").EntireRow.Insert()

# Row 212 - new record (Pehuenche)
$ws.Range("A212").Value = 4
$ws.Range("B212").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C212").Value = "Los Lagos"
$ws.Range("D212").Value = 44508
$ws.Range("E212").Value = 10
$ws.Range("F212").Value = 100114001
$ws.Range("G212").Value = "Papa"
$ws.Range("H212").Value = "Pehuenche"
$ws.Range("I212").Value = "1a nueva(o)"
$ws.Range("J212").Value = 150
$ws.Range("K212").Value = 16000
$ws.Range("L212").Value = 16000
$ws.Range("M212").Value = 16000
$ws.Range("N212").Value = "`$/saco 25 kilos"
$ws.Range("O212").Value = "Región de La Araucanía"
$ws.Range("P212").Value = 640
$ws.Range("Q212").Value = 25
$ws.Range("R212").Value = "Hortaliza"

# Row 213 - new record (Pukará)
$ws.Range("A213").Value = 4
$ws.Range("B213").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C213").Value = "Los Lagos"
$ws.Range("D213").Value = 44508
$ws.Range("E213").Value = 10
$ws.Range("F213").Value = 100114001
$ws.Range("G213").Value = "Papa"
$ws.Range("H213").Value = "Pukará"
$ws.Range("I213").Value = "1a nueva(o)"
$ws.Range("J213").Value = 150
$ws.Range("K213").Value = 16000
$ws.Range("L213").Value = 16000
$ws.Range("M213").Value = 16000
$ws.Range("N213").Value = "`$/saco 25 kilos"
$ws.Range("O213").Value = "Región de La Araucanía"
$ws.Range("P213").Value = 640
$ws.Range("Q213").Value = 25
$ws.Range("R213").Value = "Hortaliza"

# Ensure the date column keeps the date number format used elsewhere in column D
$ws.Range("D212:D213").NumberFormat = "YYYY-MM-DD HH:MM:SS"
